$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.6647940074906367
$ws1.Range("C2").Value = 0.6225626740947076
$ws1.Range("D2").Value = 0.8370786516853933
$ws1.Range("E2").Value = 0.7140575079872205
$ws1.Range("F2").Value = 0.7831114225648214
$ws1.Range("G2").Value = 0.8261302246232585
$ws1.Range("H2").Value = 0.7034745893475851
$ws1.Range("I2").Value = 447
$ws1.Range("J2").Value = 271
$ws1.Range("K2").Value = 263
$ws1.Range("L2").Value = 87

# --- Sheet 2: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")
$ws2.Range("B2").Value = 0.7514285714285714
$ws2.Range("C2").Value = 0.4925093632958801
$ws2.Range("D2").Value = 0.5950226244343891

$ws2.Range("B3").Value = 0.6225626740947076
$ws2.Range("C3").Value = 0.8370786516853933
$ws2.Range("D3").Value = 0.7140575079872205

$ws2.Range("B4").Value = 0.6647940074906367
$ws2.Range("C4").Value = 0.6647940074906367
$ws2.Range("D4").Value = 0.6647940074906367
$ws2.Range("E4").Value = 0.6647940074906367

$ws2.Range("B5").Value = 0.6869956227616395
$ws2.Range("C5").Value = 0.6647940074906367
$ws2.Range("D5").Value = 0.6545400662108047

$ws2.Range("B6").Value = 0.6869956227616396
$ws2.Range("C6").Value = 0.6647940074906367
$ws2.Range("D6").Value = 0.6545400662108047

# --- Sheet 3: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 263
$ws3.Range("C2").Value = 271

$ws3.Range("B3").Value = 87
$ws3.Range("C3").Value = 447
